# EPBDS-7870 SmartRules: CharRange is not supported in Smart Rules
#
# Adds a new "getCharValueRule5" rules table (mirrors the existing
# getCharValueRule1 table in rows 4-18) and its matching test table
# (mirrors the existing getCharValueRule2Test table in rows 57-69,
# minus the rule7/M row) further down the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# ---------------------------------------------------------------
# 1. New rules table "Rules String getCharValueRule5(char value)"
#    rows 119-132, copying formatting from the analogous rows 4-18
#    (the getCharValueRule / getCharValueRule1 table).
# ---------------------------------------------------------------

$srcRows = 4..18
$dstRows = 119..132

for ($i = 0; $i -lt $srcRows.Length; $i++) {
    $src = $srcRows[$i]
    $dst = $dstRows[$i]
    $ws.Range("C$src`:D$src").Copy()
    $ws.Range("C$dst`:D$dst").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Range("C119").Value = "Rules String getCharValueRule5(char value)"
$ws.Range("D119").Value = ""

$ws.Range("C120").Value = "C1"
$ws.Range("D120").Value = "RET1"

$ws.Range("C121").Value = "value"
$ws.Range("D121").Value = "result"

$ws.Range("C122").Value = "CharRange values"
$ws.Range("D122").Value = "String result"

$ws.Range("C123").Value = "Values"
$ws.Range("D123").Value = "Rate"

$ws.Range("C124").Value = "A"
$ws.Range("D124").Value = "rule1"

$ws.Range("C125").Value = "B..C"
$ws.Range("D125").Value = "rule2"

$ws.Range("C126").Value = "D .. E"
$ws.Range("D126").Value = "rule3"

$ws.Range("C127").Value = " F .. G "
$ws.Range("D127").Value = "rule4"

$ws.Range("C128").Value = "<=J"
$ws.Range("D128").Value = "rule5"

$ws.Range("C129").Value = "<= L"
$ws.Range("D129").Value = "rule6"

$ws.Range("C130").Value = ">= Y"
$ws.Range("D130").Value = "rule8"

$ws.Range("C131").Value = "> W"
$ws.Range("D131").Value = "rule9"

$ws.Range("C132").Value = ">= U"
$ws.Range("D132").Value = "rule10"

$ws.Range("C119:D119").Merge()

# ---------------------------------------------------------------
# 2. New test table "Test getCharValueRule5 getCharValueRule5Test"
#    rows 137-148, copying formatting from the analogous rows
#    57-69 (getCharValueRule2 / getCharValueRule1 test table),
#    skipping the old "M"/"rule7" row.
# ---------------------------------------------------------------

$srcRows2 = @(57,58,59,60,61,62,63,64,65,67,68,69)
$dstRows2 = 137..148

for ($i = 0; $i -lt $srcRows2.Length; $i++) {
    $src = $srcRows2[$i]
    $dst = $dstRows2[$i]
    $ws.Range("C$src`:D$src").Copy()
    $ws.Range("C$dst`:D$dst").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

$ws.Range("C137").Value = "Test getCharValueRule5 getCharValueRule5Test"
$ws.Range("D137").Value = ""

$ws.Range("C138").Value = "value"
$ws.Range("D138").Value = "_res_"

$ws.Range("C139").Value = "Value"
$ws.Range("D139").Value = "Result"

$ws.Range("C140").Value = "A"
$ws.Range("D140").Value = "rule1"

$ws.Range("C141").Value = "C"
$ws.Range("D141").Value = "rule2"

$ws.Range("C142").Value = "E"
$ws.Range("D142").Value = "rule3"

$ws.Range("C143").Value = "G"
$ws.Range("D143").Value = "rule4"

$ws.Range("C144").Value = "I"
$ws.Range("D144").Value = "rule5"

$ws.Range("C145").Value = "K"
$ws.Range("D145").Value = "rule6"

$ws.Range("C146").Value = "Z"
$ws.Range("D146").Value = "rule8"

$ws.Range("C147").Value = "X"
$ws.Range("D147").Value = "rule9"

$ws.Range("C148").Value = "V"
$ws.Range("D148").Value = "rule10"

$ws.Range("C137:D137").Merge()

# ---------------------------------------------------------------
# 3. Selection / view bookkeeping to match the saved workbook state.
# ---------------------------------------------------------------

$ws.Range("E134").Select()
